$wb = $excel.ActiveWorkbook

# --- DatosMotor: update PRC016 -> PRC017 family (edited first so new shared
#     strings land before the DatosCuenta ones, matching upstream ordering) ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "PRC017"
$wsMotor.Range("B2").Value = "ABC12SPRC017"
$wsMotor.Range("C2").Value = "ZAZ123SPRC017"
$wsMotor.Range("A2:C2").Select() | Out-Null

# --- DatosCuenta: new name/lastname + incremented numbers ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokeDoceName"
$wsCuenta.Range("B2").Value = "SmokeDoceLastName"
$wsCuenta.Range("C2").Value = 21546911
$wsCuenta.Range("D2").Value = 144

# --- DatosHogar: increment value, move selection ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 612
$wsHogar.Range("A2").Select() | Out-Null

# --- DatosAP: increment value ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21840817

# --- Active tab moves back to DatosCuenta (first sheet) and becomes the
#     selected/visible tab; selection there moves to H9 ---
$wsCuenta.Activate() | Out-Null
$wsCuenta.Range("H9").Select() | Out-Null
